$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-24 Monday" "2025-11-25 Tuesday"

Replace-Text "399÷5=" "753÷5="
Replace-Text "654÷6=" "862÷9="
Replace-Text "215÷4=" "692÷9="
Replace-Text "219÷2=" "801÷2="
Replace-Text "256÷5=" "544÷2="
Replace-Text "497÷8=" "530÷4="
Replace-Text "876÷9=" "477÷5="
Replace-Text "777÷6=" "435÷7="
Replace-Text "385÷4=" "318÷9="
Replace-Text "561÷2=" "951÷6="
Replace-Text "348÷4=" "554÷6="
Replace-Text "342÷2=" "800÷4="
Replace-Text "774÷3=" "617÷8="
Replace-Text "511÷9=" "244÷8="
Replace-Text "552÷8=" "969÷8="
Replace-Text "464÷2=" "714÷7="
Replace-Text "624÷8=" "546÷7="
Replace-Text "564÷4=" "794÷3="
Replace-Text "372÷7=" "866÷9="
Replace-Text "972÷9=" "188÷4="
Replace-Text "679÷7=" "984÷4="
Replace-Text "362÷4=" "749÷7="
Replace-Text "924÷5=" "853÷3="
Replace-Text "732÷2=" "471÷8="
Replace-Text "410÷4=" "478÷2="
